$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9007083773612976
$ws.Range("B1").Value = 2.131555795669556
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.848961234092712
$ws.Range("E1").Value = 1.11254870891571
